$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Fats Waller's birthdate (row 11): month/day/year
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 21
$ws.Range("E11").Value = 1904

# Remove the last two rows (Test/Joe and John/Testies)
$ws.Rows("15:16").Delete()

# Move selection to A16
$ws.Range("A16").Select()
